$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 634.46155
$ws.Range("J8").Value = 1600
$ws.Range("L8").Value = 4800
$ws.Range("N8").Value = -5078

$ws.Range("H62").Value = 3254.9092
$ws.Range("I62").Value = 3254.9092
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3254.9092
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2630.9092
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3254.9092
$ws.Range("I65").Value = 3254.9092
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16274.546
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13154.546
$ws.Range("N65").ClearContents()

$ws.Range("H113").Value = 5360
$ws.Range("J113").Value = 7333.3335
$ws.Range("L113").Value = 7333.3335
$ws.Range("N113").Value = -13841.3335

$ws.Range("H137").Value = 556982.25
$ws.Range("I137").Value = 1607.55
$ws.Range("J137").Value = 927232.0600000001
$ws.Range("K137").Value = 4822.65
$ws.Range("L137").Value = 2781696.18
$ws.Range("M137").Value = -2272.65
$ws.Range("N137").Value = -2786796.18

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 3739
$ws.Range("I141").Value = 3739
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 11217
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -6037
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2316.6667
$ws.Range("I12").Value = 1725
$ws.Range("J12").Value = 3500
$ws.Range("K12").Value = 1725
$ws.Range("L12").Value = 3500
$ws.Range("M12").Value = -1552
$ws.Range("N12").Value = -3846

$ws.Range("H61").Value = 3295.3809
$ws.Range("I61").Value = 1843.4445
$ws.Range("J61").Value = 4384.3335
$ws.Range("K61").Value = 1843.4445
$ws.Range("L61").Value = 4384.3335
$ws.Range("M61").Value = -1631.4445
$ws.Range("N61").Value = -4808.3335

$ws.Range("H63").Value = 3520
$ws.Range("I63").Value = 2540
$ws.Range("K63").Value = 2540
$ws.Range("M63").Value = -1854

$ws.Range("H66").Value = 3520
$ws.Range("I66").Value = 2540
$ws.Range("K66").Value = 12700
$ws.Range("M66").Value = -9268

$ws.Range("H74").Value = 34401.3
$ws.Range("I74").Value = 50812.85
$ws.Range("J74").Value = 1578.2
$ws.Range("K74").Value = 50812.85
$ws.Range("L74").Value = 1578.2
$ws.Range("M74").Value = -49938.85
$ws.Range("N74").Value = -3326.2

$ws.Range("H77").Value = 34401.3
$ws.Range("I77").Value = 50812.85
$ws.Range("J77").Value = 1578.2
$ws.Range("K77").Value = 254064.25
$ws.Range("L77").Value = 7891
$ws.Range("M77").Value = -249696.25
$ws.Range("N77").Value = -16627

$ws.Range("H132").Value = 2675.35
$ws.Range("I132").Value = 2154.6428
$ws.Range("J132").Value = 3890.3333
$ws.Range("K132").Value = 6463.928400000001
$ws.Range("L132").Value = 11670.9999
$ws.Range("M132").Value = -3933.928400000001
$ws.Range("N132").Value = -16730.9999

$ws.Range("H136").Value = 3295.3809
$ws.Range("I136").Value = 1843.4445
$ws.Range("J136").Value = 4384.3335
$ws.Range("K136").Value = 5530.333500000001
$ws.Range("L136").Value = 13153.0005
$ws.Range("M136").Value = -2980.333500000001
$ws.Range("N136").Value = -18253.0005

$ws.Range("H137").Value = 62200
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 62200
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 62200
$ws.Range("N137").Value = -72400
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8724.368
$ws.Range("I134").Value = 10371.917
$ws.Range("J134").Value = 5900
$ws.Range("K134").Value = 31115.751
$ws.Range("L134").Value = 17700
$ws.Range("M134").Value = -28580.751
$ws.Range("N134").Value = -22770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7463702
$ws.Range("I31").Value = 713.5
$ws.Range("J31").Value = 11112274
$ws.Range("K31").Value = 713.5
$ws.Range("L31").Value = 11112274
$ws.Range("M31").Value = -418.5
$ws.Range("N31").Value = -11112864

$ws.Range("H34").Value = 7463702
$ws.Range("I34").Value = 713.5
$ws.Range("J34").Value = 11112274
$ws.Range("K34").Value = 713.5
$ws.Range("L34").Value = 11112274
$ws.Range("M34").Value = -511.5
$ws.Range("N34").Value = -11112678

$ws.Range("H59").Value = 16033.833
$ws.Range("I59").Value = 4204
$ws.Range("J59").Value = 18399.8
$ws.Range("K59").Value = 4204
$ws.Range("L59").Value = 18399.8
$ws.Range("M59").Value = -3059
$ws.Range("N59").Value = -20689.8

$ws.Range("H68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21498

$ws.Range("H71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -67488

$ws.Range("H74").Value = 25665.285
$ws.Range("J74").Value = 25665.285
$ws.Range("L74").Value = 25665.285
$ws.Range("N74").Value = -27413.285

$ws.Range("H77").Value = 25665.285
$ws.Range("J77").Value = 25665.285
$ws.Range("L77").Value = 76995.855
$ws.Range("N77").Value = -85731.855

$ws.Range("H99").Value = 1913.95
$ws.Range("I99").Value = 1949.3125
$ws.Range("K99").Value = 1949.3125
$ws.Range("M99").Value = -451.3125

$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -29868

$ws.Range("H104").Value = 28270.6
$ws.Range("J104").Value = 28270.6
$ws.Range("L104").Value = 28270.6
$ws.Range("N104").Value = -33512.6

$ws.Range("H120").Value = 27500
$ws.Range("J120").Value = 27500
$ws.Range("L120").Value = 27500
$ws.Range("N120").Value = -34758

$ws.Range("H126").Value = 1913.95
$ws.Range("I126").Value = 1949.3125
$ws.Range("K126").Value = 5847.9375
$ws.Range("M126").Value = -3377.9375

$ws.Range("H134").Value = 2326.4167
$ws.Range("I134").Value = 2196.8572
$ws.Range("J134").Value = 3233.3333
$ws.Range("K134").Value = 6590.571599999999
$ws.Range("L134").Value = 9699.999899999999
$ws.Range("M134").Value = -4055.571599999999
$ws.Range("N134").Value = -14769.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1815.95
$ws.Range("I68").Value = 1034.85
$ws.Range("J68").Value = 2206.5
$ws.Range("K68").Value = 3104.55
$ws.Range("L68").Value = 6619.5
$ws.Range("M68").Value = -2293.55
$ws.Range("N68").Value = -8241.5

$ws.Range("H71").Value = 1815.95
$ws.Range("I71").Value = 1034.85
$ws.Range("J71").Value = 2206.5
$ws.Range("K71").Value = 9313.65
$ws.Range("L71").Value = 19858.5
$ws.Range("M71").Value = -5257.65
$ws.Range("N71").Value = -27970.5

$ws.Range("H107").Value = 723266.9
$ws.Range("I107").Value = 1054.0769
$ws.Range("J107").Value = 923027.9
$ws.Range("K107").Value = 3162.2307
$ws.Range("L107").Value = 2769083.7
$ws.Range("M107").Value = -1242.2307
$ws.Range("N107").Value = -2772923.7

$ws.Range("H136").Value = 10000
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3230.0688
$ws.Range("I132").Value = 4486
$ws.Range("J132").Value = 2664.9
$ws.Range("K132").Value = 13458
$ws.Range("L132").Value = 7994.700000000001
$ws.Range("M132").Value = -10928
$ws.Range("N132").Value = -13054.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 362.54285
$ws.Range("I55").Value = 167.69565
$ws.Range("J55").Value = 736
$ws.Range("K55").Value = 167.69565
$ws.Range("L55").Value = 736
$ws.Range("M55").Value = 5.304349999999999
$ws.Range("N55").Value = -1082

$ws.Range("H132").Value = 3050.1187
$ws.Range("I132").Value = 2524.842
$ws.Range("J132").Value = 4000.6191
$ws.Range("K132").Value = 7574.526
$ws.Range("L132").Value = 12001.8573
$ws.Range("M132").Value = -5044.526
$ws.Range("N132").Value = -17061.8573

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H122").Value = 49742.12
$ws.Range("I122").Value = 1754.1428
$ws.Range("J122").Value = 110817.73
$ws.Range("K122").Value = 5262.428400000001
$ws.Range("L122").Value = 332453.19
$ws.Range("M122").Value = -2812.428400000001
$ws.Range("N122").Value = -337353.19

$ws.Range("H126").Value = 90910880
$ws.Range("I126").Value = 1364.8334
$ws.Range("K126").Value = 4094.5002
$ws.Range("M126").Value = -1624.5002

$ws.Range("H136").Value = 2499.325
$ws.Range("I136").Value = 2760.7693
$ws.Range("K136").Value = 8282.3079
$ws.Range("M136").Value = -5732.3079
